# Add a new "CarNameAndPrice" worksheet between "FindCarTest" and
# "UserRegTest", populate it with carbrand/browserType/runmode test data
# (introducing a new "firefox" shared string), tweak the FindCarTest
# sheet's selection/column width, and make the new sheet the active tab.

$wb = $excel.ActiveWorkbook

$findCarSheet = $wb.Worksheets.Item("FindCarTest")

# FindCarTest: selection becomes A1:C4 (no longer the tab-selected sheet,
# that moves to the new sheet below), and a width is set on column D.
$findCarSheet.Range("A1:C4").Select() | Out-Null
$findCarSheet.Columns.Item(4).ColumnWidth = 14

# Insert the new sheet right after FindCarTest.
$newSheet = $wb.Worksheets.Add($null, $findCarSheet)
$newSheet.Name = "CarNameAndPrice"

$newSheet.Cells.Item(1,1).Value = "carbrand"
$newSheet.Cells.Item(1,2).Value = "browserType"
$newSheet.Cells.Item(1,3).Value = "runmode"

$newSheet.Cells.Item(2,1).Value = "bmw"
$newSheet.Cells.Item(2,2).Value = "chrome"
$newSheet.Cells.Item(2,3).Value = "y"

$newSheet.Cells.Item(3,1).Value = "mg"
$newSheet.Cells.Item(3,2).Value = "firefox"
$newSheet.Cells.Item(3,3).Value = "y"

$newSheet.Cells.Item(4,1).Value = "toyota"
$newSheet.Cells.Item(4,2).Value = "chrome"
$newSheet.Cells.Item(4,3).Value = "y"

$newSheet.Range("C6").Select() | Out-Null
$newSheet.Activate() | Out-Null
